# Generate Report for Handback
# Update the localization-status workbook with handback failure information.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Row 3 on the Overview sheet corresponds to file
# 23ec9865-7ccb-4d28-9641-802e4b3db46d.md ; its Status changed from
# "Ready for handoff" to "Handback transform failed".
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# The per-language sheets share the same Status text for that file (column C).
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Populate the Error Detail column (P) for that row with the handback
# transform failure explanation, and widen the column to fit the text
# (raw OOXML column width of 40 == Excel ColumnWidth ~39.1666667 chars).
$wsZhCn.Range("P3").Value = "Handback file name: 4sooc3xj.cef is different with handoff file name: 23ec9865-7ccb-4d28-9641-802e4b3db46d.9240f9aba8822260a902b07ed63700e29df763a3.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666667

$wsDeDe.Range("P3").Value = "Handback file name: 4sooc3xj.cef is different with handoff file name: 23ec9865-7ccb-4d28-9641-802e4b3db46d.9240f9aba8822260a902b07ed63700e29df763a3.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666667
